$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest cryptos data refresh: update Price (D) / Volume(1h) (E) figures,
# and swap the PancakeSwap / Decentraland rows (45-46) to match new ranking order.

$ws.Range("D2").Value = "27.573.54"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.749.53"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.73"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4591"
$ws.Range("E7").Value = "  +9.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3569"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07481"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.02"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.092"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.74"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.997"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.085"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "1.751.19"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.34"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06424"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.797"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").Value = "27.637.54"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.108"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.53"
$ws.Range("E26").Value = "  +4.94%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "1.950.04"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.93"
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.053"
$ws.Range("E31").Value = "  -6.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09210"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.668"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.529"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02294"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06039"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2087"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.972"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6297"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.377"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.748"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.17"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.716"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5887"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.44"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.938"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.93"
$ws.Range("E51").Value = "  -1.93%  "
